$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price(D) and Volume(E) columns are treated as text so that
# numeric-looking strings (e.g. "181.70", "0.0000119") are not silently
# converted into numbers/scientific notation, and formatting like
# trailing zeros or thousand-dot separators is preserved.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.854.73"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -3.70%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.344.50"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.72%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.42"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -3.23%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.70"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -5.20%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.44%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -3.12%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.66"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.51%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -4.15%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.926.31"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.79%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.71%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.21"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -4.95%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "66.891.70"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.74%  "

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.20%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.333.81"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.79%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "436.23"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -3.31%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.64"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.35%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.66%  "

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.30%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.89"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.98%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.30%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.25%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000119"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.21%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.21%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.05"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -4.87%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.14%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.97"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.26%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.89"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.74%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -5.75%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.999"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.00%  "

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.48%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.82"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.71%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.07%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "161.51"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.25%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "27.98"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.98%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -4.63%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.836.87"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +3.76%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -3.08%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.47%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.31"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.69%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.49%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.17"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.10%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "24.69"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.91%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -6.94%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "325.08"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -5.48%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -4.17%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.990"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.99%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "31.23"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -5.09%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.17"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.66%  "
